$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 258.85715
$ws.Range("I11").Value = 258.85715
$ws.Range("K11").Value = 258.85715
$ws.Range("M11").Value = -118.85715
$ws.Range("H38").Value = 124.8
$ws.Range("I38").Value = 27.555555
$ws.Range("J38").Value = 1000
$ws.Range("K38").Value = 82.66666499999999
$ws.Range("L38").Value = 3000
$ws.Range("M38").Value = 289.333335
$ws.Range("N38").Value = -3744
$ws.Range("H69").Value = 3964.5417
$ws.Range("I69").Value = 3983.3333
$ws.Range("J69").Value = 3953.2666
$ws.Range("K69").Value = 11949.9999
$ws.Range("L69").Value = 11859.7998
$ws.Range("M69").Value = -11075.9999
$ws.Range("N69").Value = -13607.7998
$ws.Range("H70").Value = 1877.2325
$ws.Range("I70").Value = 6385
$ws.Range("J70").Value = 1414.8975
$ws.Range("K70").Value = 19155
$ws.Range("L70").Value = 4244.6925
$ws.Range("M70").Value = -18885
$ws.Range("N70").Value = -4784.6925
$ws.Range("H72").Value = 3964.5417
$ws.Range("I72").Value = 3983.3333
$ws.Range("J72").Value = 3953.2666
$ws.Range("K72").Value = 35849.9997
$ws.Range("L72").Value = 35579.3994
$ws.Range("M72").Value = -31481.9997
$ws.Range("N72").Value = -44315.3994
$ws.Range("H73").Value = 1877.2325
$ws.Range("I73").Value = 6385
$ws.Range("J73").Value = 1414.8975
$ws.Range("K73").Value = 19155
$ws.Range("L73").Value = 4244.6925
$ws.Range("M73").Value = -18219
$ws.Range("N73").Value = -6116.6925
$ws.Range("H74").Value = 3000.5715
$ws.Range("I74").Value = 2652
$ws.Range("J74").Value = 3140
$ws.Range("K74").Value = 2652
$ws.Range("L74").Value = 3140
$ws.Range("M74").Value = -1716
$ws.Range("N74").Value = -5012
$ws.Range("H77").Value = 3000.5715
$ws.Range("I77").Value = 2652
$ws.Range("J77").Value = 3140
$ws.Range("K77").Value = 13260
$ws.Range("L77").Value = 15700
$ws.Range("M77").Value = -8580
$ws.Range("N77").Value = -25060
$ws.Range("H80").Value = 568.80554
$ws.Range("I80").Value = 326.47058
$ws.Range("J80").Value = 785.6316
$ws.Range("K80").Value = 979.41174
$ws.Range("L80").Value = 2356.8948
$ws.Range("M80").Value = 18.58825999999999
$ws.Range("N80").Value = -4352.8948
$ws.Range("H83").Value = 568.80554
$ws.Range("I83").Value = 326.47058
$ws.Range("J83").Value = 785.6316
$ws.Range("K83").Value = 2938.23522
$ws.Range("L83").Value = 7070.6844
$ws.Range("M83").Value = 2053.76478
$ws.Range("N83").Value = -17054.6844
$ws.Range("H96").Value = 406.29413
$ws.Range("I96").Value = 344.46155
$ws.Range("J96").Value = 607.25
$ws.Range("K96").Value = 1033.38465
$ws.Range("L96").Value = 1821.75
$ws.Range("M96").Value = 339.61535
$ws.Range("N96").Value = -4567.75
$ws.Range("H100").Value = 47393.184
$ws.Range("I100").Value = 56818.332
$ws.Range("J100").Value = 4980
$ws.Range("K100").Value = 56818.332
$ws.Range("L100").Value = 4980
$ws.Range("M100").Value = -56277.332
$ws.Range("N100").Value = -6062
$ws.Range("H103").Value = 402
$ws.Range("I103").Value = 401.5
$ws.Range("J103").Value = 405
$ws.Range("K103").Value = 1204.5
$ws.Range("L103").Value = 1215
$ws.Range("M103").Value = -618.5
$ws.Range("N103").Value = -2387
$ws.Range("H121").Value = 1278.5714
$ws.Range("I121").Value = 575
$ws.Range("J121").Value = 1560
$ws.Range("K121").Value = 1725
$ws.Range("L121").Value = 4680
$ws.Range("M121").Value = 22
$ws.Range("N121").Value = -8174

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 57.625
$ws.Range("I4").Value = 57.625
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 57.625
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 58.375
$ws.Range("N4").ClearContents()
$ws.Range("H74").Value = 4703.9644
$ws.Range("I74").Value = 1141.7222
$ws.Range("J74").Value = 11116
$ws.Range("K74").Value = 1141.7222
$ws.Range("L74").Value = 11116
$ws.Range("M74").Value = -267.7221999999999
$ws.Range("N74").Value = -12864
$ws.Range("H77").Value = 4703.9644
$ws.Range("I77").Value = 1141.7222
$ws.Range("J77").Value = 11116
$ws.Range("K77").Value = 5708.611
$ws.Range("L77").Value = 55580
$ws.Range("M77").Value = -1340.611
$ws.Range("N77").Value = -64316
$ws.Range("H88").Value = 4169
$ws.Range("I88").Value = 3500
$ws.Range("J88").Value = 4503.5
$ws.Range("K88").Value = 3500
$ws.Range("L88").Value = 4503.5
$ws.Range("M88").Value = -3094
$ws.Range("N88").Value = -5315.5
$ws.Range("H91").Value = 4169
$ws.Range("I91").Value = 3500
$ws.Range("J91").Value = 4503.5
$ws.Range("K91").Value = 3500
$ws.Range("L91").Value = 4503.5
$ws.Range("M91").Value = -2096
$ws.Range("N91").Value = -7311.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 366.1111
$ws.Range("I64").Value = 208.66667
$ws.Range("K64").Value = 208.66667
$ws.Range("M64").Value = 16.33332999999999
$ws.Range("H67").Value = 366.1111
$ws.Range("I67").Value = 208.66667
$ws.Range("K67").Value = 208.66667
$ws.Range("M67").Value = 571.3333299999999
$ws.Range("H86").Value = 4659.8335
$ws.Range("I86").Value = 4326.3335
$ws.Range("J86").Value = 4993.3335
$ws.Range("K86").Value = 4326.3335
$ws.Range("L86").Value = 4993.3335
$ws.Range("M86").Value = -3203.3335
$ws.Range("N86").Value = -7239.3335
$ws.Range("H89").Value = 4659.8335
$ws.Range("I89").Value = 4326.3335
$ws.Range("J89").Value = 4993.3335
$ws.Range("K89").Value = 21631.6675
$ws.Range("L89").Value = 24966.6675
$ws.Range("M89").Value = -16015.6675
$ws.Range("N89").Value = -36198.6675

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 7995
$ws.Range("J41").Value = 7995
$ws.Range("L41").Value = 7995
$ws.Range("N41").Value = -8851
$ws.Range("H58").Value = 2841
$ws.Range("I58").Value = 1040
$ws.Range("J58").Value = 5842.6665
$ws.Range("K58").Value = 1040
$ws.Range("L58").Value = 5842.6665
$ws.Range("M58").Value = -837
$ws.Range("N58").Value = -6248.6665
$ws.Range("H62").Value = 3779.1667
$ws.Range("I62").Value = 3985
$ws.Range("J62").Value = 2750
$ws.Range("K62").Value = 3985
$ws.Range("L62").Value = 2750
$ws.Range("M62").Value = -3361
$ws.Range("N62").Value = -3998
$ws.Range("H65").Value = 3779.1667
$ws.Range("I65").Value = 3985
$ws.Range("J65").Value = 2750
$ws.Range("K65").Value = 19925
$ws.Range("L65").Value = 13750
$ws.Range("M65").Value = -16805
$ws.Range("N65").Value = -19990
$ws.Range("H136").Value = 2841
$ws.Range("I136").Value = 1040
$ws.Range("J136").Value = 5842.6665
$ws.Range("K136").Value = 3120
$ws.Range("L136").Value = 17527.9995
$ws.Range("M136").Value = -570
$ws.Range("N136").Value = -22627.9995

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2400
$ws.Range("I39").Value = 500
$ws.Range("J39").Value = 3666.6667
$ws.Range("K39").Value = 1500
$ws.Range("L39").Value = 11000.0001
$ws.Range("M39").Value = -1206
$ws.Range("N39").Value = -11588.0001
$ws.Range("H131").Value = 654.5179000000001
$ws.Range("I131").Value = 331.36
$ws.Range("J131").Value = 915.129
$ws.Range("K131").Value = 994.08
$ws.Range("L131").Value = 2745.387
$ws.Range("M131").Value = 4045.92
$ws.Range("N131").Value = -12825.387

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 152.5
$ws.Range("I13").Value = 152.5
$ws.Range("K13").Value = 152.5
$ws.Range("M13").Value = -13.5
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H64").Value = 29271
$ws.Range("J64").Value = 29271
$ws.Range("L64").Value = 29271
$ws.Range("N64").Value = -29767
$ws.Range("H67").Value = 29271
$ws.Range("J67").Value = 29271
$ws.Range("L67").Value = 29271
$ws.Range("N67").Value = -30987

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2756.8462
$ws.Range("I22").Value = 1311.4286
$ws.Range("J22").Value = 4443.1665
$ws.Range("K22").Value = 1311.4286
$ws.Range("L22").Value = 4443.1665
$ws.Range("M22").Value = -1016.4286
$ws.Range("N22").Value = -5033.1665
$ws.Range("H27").Value = 2756.8462
$ws.Range("I27").Value = 1311.4286
$ws.Range("J27").Value = 4443.1665
$ws.Range("K27").Value = 1311.4286
$ws.Range("L27").Value = 4443.1665
$ws.Range("M27").Value = -1204.4286
$ws.Range("N27").Value = -4657.1665
$ws.Range("H46").Value = 2254.7
$ws.Range("I46").Value = 3659.6667
$ws.Range("J46").Value = 1652.5714
$ws.Range("K46").Value = 3659.6667
$ws.Range("L46").Value = 1652.5714
$ws.Range("M46").Value = -3471.6667
$ws.Range("N46").Value = -2028.5714
$ws.Range("H55").Value = 116.28571
$ws.Range("I55").Value = 112.4
$ws.Range("J55").Value = 126
$ws.Range("K55").Value = 112.4
$ws.Range("L55").Value = 126
$ws.Range("M55").Value = 60.59999999999999
$ws.Range("N55").Value = -472
$ws.Range("H136").Value = 5659.407
$ws.Range("I136").Value = 1992.909
$ws.Range("J136").Value = 21792
$ws.Range("K136").Value = 5978.727000000001
$ws.Range("L136").Value = 65376
$ws.Range("M136").Value = -3428.727000000001
$ws.Range("N136").Value = -70476

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2440.2
$ws.Range("I96").Value = 2100.5
$ws.Range("J96").Value = 2666.6667
$ws.Range("K96").Value = 2666.6667
$ws.Range("L96").Value = 2666.6667
$ws.Range("M96").Value = -727.5
$ws.Range("N96").Value = -5412.6667
